$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "<the>"
$ws.Range("C2").Value = 31

$ws.Range("C3").Value = 32

$ws.Range("C4").Value = 28

$ws.Range("B5").Value = "<yilo>"
$ws.Range("C5").Value = 30

$ws.Range("C6").Value = 27

$ws.Range("C7").Value = 31

$ws.Range("B8").Value = "<four>"
$ws.Range("C8").Value = 24

$ws.Range("C9").Value = 35

$ws.Range("C10").Value = 30

$ws.Range("B11").Value = "<in>"
$ws.Range("C11").Value = 26

$ws.Range("C12").Value = 27

$ws.Range("C14").Value = 35

$ws.Range("B15").Value = "<in>"
$ws.Range("C15").Value = 30

$ws.Range("C16").Value = 33

$ws.Range("B17").Value = "<so>"
$ws.Range("C17").Value = 34
